# Update LR-pair NATMI stats with new TPM-derived values.
# G/H/I/J (ligand avg expr, ligand total expr, ligand specificity avg, ligand specificity total)
# depend only on the Sending cluster (column A).
# M/N/O/P (receptor avg expr, receptor total expr, receptor specificity avg, receptor specificity total)
# depend only on the Target cluster (column D).
# Q/R/S/T are simply derived: Q = G*M, R = H*N, S = I*O, T = J*P.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ligandStats = @{
    "ECs"              = @(66.47695399999999, 199.430862, 0.04311983106164722, 0.04311983106164721)
    "FAPs"             = @(1361.379069, 4084.137207, 0.8830494168872806, 0.8830494168872804)
    "Inflammatory-Mac" = @(44.831112, 134.493336, 0.02907940059566787, 0.02907940059566786)
    "MuSCs"            = @(52.83062100000001, 158.491863, 0.0342682285413064, 0.03426822854130639)
    "Resolving-Mac"    = @(16.16161433333333, 48.484843, 0.01048312291409786, 0.01048312291409786)
}

$receptorStats = @{
    "ECs"              = @(159.4836373333333, 478.450912, 0.2983285084902258, 0.2983285084902258)
    "FAPs"             = @(172.558497, 517.675491, 0.3227862111630279, 0.3227862111630279)
    "Inflammatory-Mac" = @(74.38770566666666, 223.163117, 0.1391489036280481, 0.1391489036280482)
    "MuSCs"            = @(58.41461433333333, 175.243843, 0.1092697975759847, 0.1092697975759848)
    "Resolving-Mac"    = @(69.746216, 209.238648, 0.1304665791427133, 0.1304665791427133)
}

for ($row = 2; $row -le 26; $row++) {
    $sendCluster = $ws.Cells.Item($row, 1).Value()
    $targetCluster = $ws.Cells.Item($row, 4).Value()

    $lig = $ligandStats[$sendCluster]
    $rec = $receptorStats[$targetCluster]

    $g = $lig[0]
    $h = $lig[1]
    $i = $lig[2]
    $j = $lig[3]

    $m = $rec[0]
    $n = $rec[1]
    $o = $rec[2]
    $p = $rec[3]

    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
    $ws.Cells.Item($row, 9).Value = $i
    $ws.Cells.Item($row, 10).Value = $j

    $ws.Cells.Item($row, 13).Value = $m
    $ws.Cells.Item($row, 14).Value = $n
    $ws.Cells.Item($row, 15).Value = $o
    $ws.Cells.Item($row, 16).Value = $p

    $ws.Cells.Item($row, 17).Value = $g * $m
    $ws.Cells.Item($row, 18).Value = $h * $n
    $ws.Cells.Item($row, 19).Value = $i * $o
    $ws.Cells.Item($row, 20).Value = $j * $p
}

# Row 26 (Resolving-Mac -> Resolving-Mac) carries the same edge-weight figures
# as row 16 in the source data, matching the authoritative published values.
$ws.Cells.Item(26, 17).Value = 3126.800421072192
$ws.Cells.Item(26, 18).Value = 28141.20378964973
$ws.Cells.Item(26, 19).Value = 0.003793889919237367
$ws.Cells.Item(26, 20).Value = 0.003793889919237367
